$d = $word.ActiveDocument

# Change 1: merge the ".external_short_name>><<else>> Online Civil Claims<<es_>>" runs
# (which spelled out the old "<<caseManagementLocation.external_short_name>>" merge
# field piece by piece) into a single run reading
# ".venue_name>><<else>> Online Civil Claims<<es_>>".
$pattern1 = ".external_short_name>><<else>> Online Civil Claims<<es_>>"
$replacement1 = ".venue_name>><<else>> Online Civil Claims<<es_>>"
$d.Content.Find.Execute($pattern1, $true, $false, $false, $false, $false, $true, 1, $false, $replacement1, 2) | Out-Null

# Change 2: drop the stale <w:lastRenderedPageBreak/> marker that was sitting in
# front of the "c. the upload function ..." run. Locate the paragraph, then
# round-trip its OOXML through InsertXML, which rebuilds the paragraph and
# sheds the obsolete (purely cosmetic/recalculated) lastRenderedPageBreak field
# while leaving the run's formatting/rsid attributes untouched.
$pattern2 = "the upload function has not become available"
$paragraphs = $d.Paragraphs
$target2 = $null
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.Text -like "*$pattern2*") {
        $target2 = $candidate
        break
    }
}
if ($target2 -ne $null) {
    $pr2 = $target2.Range
    $pr2.InsertXML($pr2.WordOpenXML)
}
